$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "唯一标识"
$ws.Range("B1").Value = "名字"
$ws.Range("C1").Value = "how old"

$ws.Range("C2").Select()
